$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 1756.5714
$ws.Range("I34").Value = 1756.5714
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1756.5714
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1553.5714
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 1756.5714
$ws.Range("I36").Value = 1756.5714
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1756.5714
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1041.5714
$ws.Range("N36").ClearContents()
$ws.Range("H88").Value = 8484.143
$ws.Range("I88").Value = 500
$ws.Range("K88").Value = 500
$ws.Range("M88").Value = -94
$ws.Range("H91").Value = 8484.143
$ws.Range("I91").Value = 500
$ws.Range("K91").Value = 500
$ws.Range("M91").Value = 904
$ws.Range("H132").Value = 1636.4255
$ws.Range("I132").Value = 1486.0476
$ws.Range("J132").Value = 2899.6
$ws.Range("K132").Value = 4458.142800000001
$ws.Range("L132").Value = 8698.799999999999
$ws.Range("M132").Value = -1928.142800000001
$ws.Range("N132").Value = -13758.8
$ws.Range("H135").Value = 2709.652
$ws.Range("I135").Value = 2767.889
$ws.Range("K135").Value = 24911.001
$ws.Range("M135").Value = -22376.001
$ws.Range("H137").Value = 1538.037
$ws.Range("I137").Value = 1251.8823
$ws.Range("J137").Value = 2024.5
$ws.Range("K137").Value = 3755.6469
$ws.Range("L137").Value = 6073.5
$ws.Range("M137").Value = -1205.6469
$ws.Range("N137").Value = -11173.5
$ws.Range("H141").Value = 10046
$ws.Range("I141").Value = 12847.667
$ws.Range("J141").Value = 3322
$ws.Range("K141").Value = 38543.001
$ws.Range("L141").Value = 9966
$ws.Range("M141").Value = -33363.001
$ws.Range("N141").Value = -20326

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 38470284
$ws.Range("I61").Value = 38470284
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 38470284
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -38470072
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 22729218
$ws.Range("I74").Value = 31251504
$ws.Range("J74").Value = 3125.0833
$ws.Range("K74").Value = 31251504
$ws.Range("L74").Value = 3125.0833
$ws.Range("M74").Value = -31250630
$ws.Range("N74").Value = -4873.0833
$ws.Range("H77").Value = 22729218
$ws.Range("I77").Value = 31251504
$ws.Range("J77").Value = 3125.0833
$ws.Range("K77").Value = 156257520
$ws.Range("L77").Value = 15625.4165
$ws.Range("M77").Value = -156253152
$ws.Range("N77").Value = -24361.4165
$ws.Range("H132").Value = 38477524
$ws.Range("I132").Value = 6230.048
$ws.Range("J132").Value = 200056960
$ws.Range("K132").Value = 18690.144
$ws.Range("L132").Value = 600170880
$ws.Range("M132").Value = -16160.144
$ws.Range("N132").Value = -600175940
$ws.Range("H136").Value = 38470284
$ws.Range("I136").Value = 38470284
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 115410852
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -115408302
$ws.Range("N136").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 33510.8
$ws.Range("I5").Value = 2514.6667
$ws.Range("K5").Value = 2514.6667
$ws.Range("M5").Value = -2401.6667
$ws.Range("H86").Value = 11710.909
$ws.Range("I86").Value = 11682
$ws.Range("K86").Value = 11682
$ws.Range("M86").Value = -10559
$ws.Range("H89").Value = 11710.909
$ws.Range("I89").Value = 11682
$ws.Range("K89").Value = 58410
$ws.Range("M89").Value = -52794
$ws.Range("H105").Value = 9253.883
$ws.Range("I105").Value = 9676
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 9676
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -7929
$ws.Range("N105").Value = -5994
$ws.Range("H107").Value = 866.74194
$ws.Range("I107").Value = 804.6429000000001
$ws.Range("J107").Value = 1446.3334
$ws.Range("K107").Value = 804.6429000000001
$ws.Range("L107").Value = 1446.3334
$ws.Range("M107").Value = 1115.3571
$ws.Range("N107").Value = -5286.3334
$ws.Range("H134").Value = 3326.913
$ws.Range("I134").Value = 2326.1
$ws.Range("J134").Value = 9999
$ws.Range("K134").Value = 6978.299999999999
$ws.Range("L134").Value = 29997
$ws.Range("M134").Value = -4443.299999999999
$ws.Range("N134").Value = -35067

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 150.11765
$ws.Range("I7").Value = 71.78261000000001
$ws.Range("J7").Value = 313.9091
$ws.Range("K7").Value = 71.78261000000001
$ws.Range("L7").Value = 313.9091
$ws.Range("M7").Value = 41.21738999999999
$ws.Range("N7").Value = -539.9091000000001
$ws.Range("H31").Value = 6722.636
$ws.Range("J31").Value = 7499.1665
$ws.Range("L31").Value = 7499.1665
$ws.Range("N31").Value = -8089.1665
$ws.Range("H34").Value = 6722.636
$ws.Range("J34").Value = 7499.1665
$ws.Range("L34").Value = 7499.1665
$ws.Range("N34").Value = -7903.1665
$ws.Range("H51").Value = 26077.4
$ws.Range("J51").Value = 30148.5
$ws.Range("L51").Value = 30148.5
$ws.Range("N51").Value = -31620.5
$ws.Range("H61").Value = 26077.4
$ws.Range("J61").Value = 30148.5
$ws.Range("L61").Value = 30148.5
$ws.Range("N61").Value = -30844.5
$ws.Range("H86").Value = 7773.25
$ws.Range("I86").Value = 7300.6
$ws.Range("J86").Value = 8110.857
$ws.Range("K86").Value = 7300.6
$ws.Range("L86").Value = 8110.857
$ws.Range("M86").Value = -6177.6
$ws.Range("N86").Value = -10356.857
$ws.Range("H89").Value = 7773.25
$ws.Range("I89").Value = 7300.6
$ws.Range("J89").Value = 8110.857
$ws.Range("K89").Value = 36503
$ws.Range("L89").Value = 40554.285
$ws.Range("M89").Value = -30887
$ws.Range("N89").Value = -51786.285
$ws.Range("H132").Value = 2918.2
$ws.Range("I132").Value = 2606.739
$ws.Range("K132").Value = 7820.217000000001
$ws.Range("M132").Value = -5290.217000000001
$ws.Range("H134").Value = 3195.0256
$ws.Range("I134").Value = 2620.6428
$ws.Range("J134").Value = 4657.091
$ws.Range("K134").Value = 7861.928400000001
$ws.Range("L134").Value = 13971.273
$ws.Range("M134").Value = -5326.928400000001
$ws.Range("N134").Value = -19041.273

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 182.14285
$ws.Range("I15").Value = 80.75
$ws.Range("J15").Value = 317.33334
$ws.Range("K15").Value = 242.25
$ws.Range("L15").Value = 952.0000200000001
$ws.Range("M15").Value = -102.25
$ws.Range("N15").Value = -1232.00002
$ws.Range("H32").Value = 1000000
$ws.Range("I32").Value = 1000000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3000000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2999717
$ws.Range("N32").ClearContents()
$ws.Range("H108").Value = 4850
$ws.Range("I108").Value = 4850
$ws.Range("K108").Value = 14550
$ws.Range("M108").Value = -11670
$ws.Range("H124").Value = 1894.6
$ws.Range("I124").Value = 1697.1428
$ws.Range("J124").Value = 2355.3333
$ws.Range("K124").Value = 5091.428400000001
$ws.Range("L124").Value = 7065.999899999999
$ws.Range("M124").Value = -181.4284000000007
$ws.Range("N124").Value = -16885.9999
$ws.Range("H131").Value = 3772.2
$ws.Range("I131").Value = 2025
$ws.Range("J131").Value = 4041
$ws.Range("K131").Value = 6075
$ws.Range("L131").Value = 12123
$ws.Range("M131").Value = -1035
$ws.Range("N131").Value = -22203
$ws.Range("H134").Value = 2416.2273
$ws.Range("I134").Value = 511.72223
$ws.Range("J134").Value = 10986.5
$ws.Range("K134").Value = 1535.16669
$ws.Range("L134").Value = 32959.5
$ws.Range("M134").Value = 3534.83331
$ws.Range("N134").Value = -43099.5
$ws.Range("H137").Value = 3016.5
$ws.Range("J137").Value = 3016.5
$ws.Range("L137").Value = 9049.5
$ws.Range("N137").Value = -19249.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6301.325
$ws.Range("I132").Value = 5646.9033
$ws.Range("K132").Value = 16940.7099
$ws.Range("M132").Value = -14410.7099

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7150.25
$ws.Range("I7").Value = 6458.9165
$ws.Range("K7").Value = 6458.9165
$ws.Range("M7").Value = -6346.9165
$ws.Range("H126").Value = 7150.25
$ws.Range("I126").Value = 6458.9165
$ws.Range("K126").Value = 19376.7495
$ws.Range("M126").Value = -16906.7495
$ws.Range("H132").Value = 4374.75
$ws.Range("I132").Value = 3999.5
$ws.Range("J132").Value = 4750
$ws.Range("K132").Value = 11998.5
$ws.Range("L132").Value = 14250
$ws.Range("M132").Value = -9468.5
$ws.Range("N132").Value = -19310
$ws.Range("H136").Value = 1113183
$ws.Range("I136").Value = 1334979.6
$ws.Range("J136").Value = 4199.6665
$ws.Range("K136").Value = 4004938.8
$ws.Range("L136").Value = 12598.9995
$ws.Range("M136").Value = -4002388.8
$ws.Range("N136").Value = -17698.9995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 1000
$ws.Range("I54").Value = 1000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 1000
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -480
$ws.Range("N54").ClearContents()
$ws.Range("H95").Value = 26548.5
$ws.Range("J95").Value = 26548.5
$ws.Range("L95").Value = 26548.5
$ws.Range("N95").Value = -32040.5
$ws.Range("H100").Value = 72144020
$ws.Range("I100").Value = 91818840
$ws.Range("K100").Value = 183637680
$ws.Range("M100").Value = -183637139
$ws.Range("H113").Value = 1119.9
$ws.Range("I113").Value = 1147.6316
$ws.Range("J113").Value = 593
$ws.Range("K113").Value = 3442.8948
$ws.Range("L113").Value = 1779
$ws.Range("M113").Value = -1272.8948
$ws.Range("N113").Value = -6119
$ws.Range("H125").Value = 69999
$ws.Range("J125").Value = 69999
$ws.Range("L125").Value = 69999
$ws.Range("N125").Value = -79839
$ws.Range("H132").Value = 3329.3333
$ws.Range("J132").Value = 4852.857
$ws.Range("L132").Value = 14558.571
$ws.Range("N132").Value = -19618.571
$ws.Range("H136").Value = 2750
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 3333.3333
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 9999.999899999999
$ws.Range("M136").Value = -450
$ws.Range("N136").Value = -15099.9999
